# Update "want-to-go" count (column F) values for the matching events
# on sheet "展览" (Exhibition) and sheet "全部类型" (All types).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 2157
$ws1.Range("F4").Value = 40
$ws1.Range("F5").Value = 11355
$ws1.Range("F7").Value = 317
$ws1.Range("F8").Value = 218
$ws1.Range("F9").Value = 11297
$ws1.Range("F11").Value = 1153
$ws1.Range("F12").Value = 68
$ws1.Range("F13").Value = 1740
$ws1.Range("F14").Value = 5636
$ws1.Range("F15").Value = 104
$ws1.Range("F16").Value = 3473
$ws1.Range("F18").Value = 10

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 2157
$ws4.Range("F5").Value = 40
$ws4.Range("F7").Value = 11355
$ws4.Range("F9").Value = 317
$ws4.Range("F10").Value = 218
$ws4.Range("F11").Value = 11297
$ws4.Range("F13").Value = 1153
$ws4.Range("F14").Value = 68
$ws4.Range("F15").Value = 1740
$ws4.Range("F17").Value = 5636
$ws4.Range("F18").Value = 104
$ws4.Range("F19").Value = 3473
$ws4.Range("F21").Value = 10
